$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column headers (row 1) ---
$ws.Range("F1").Value = "Falling Edge"
$ws.Range("G1").Value = "Rising Edge"
$ws.Range("C1").Value = "Rising Edge Error (deg)"
$ws.Range("B1").Value = "Falling Edge Error (deg)"

# --- Column B: duplicate of column A (Falling Edge values) ---
$ws.Range("B2").Value = 4.3
$ws.Range("B3").Value = 3.2
$ws.Range("B4").Value = -0.1
$ws.Range("B5").Value = 2.5
$ws.Range("B6").Value = 1.8
$ws.Range("B7").Value = -2.5
$ws.Range("B8").Value = 1.3
$ws.Range("B9").Value = 2.6
$ws.Range("B10").Value = 1.8
$ws.Range("B11").Value = 1.4

# --- Summary statistics block (E:G, rows 1-3) ---
$ws.Range("E2").Value = "Mean (deg)"
$ws.Range("E3").Value = "Standard Deviation (deg)"

$ws.Range("F2").Formula = "=AVERAGE(B2:B11)"
$ws.Range("G2").Formula = "=AVERAGE(C2:C11)"
$ws.Range("F3").Formula = "=STDEV(B2:B11)"
$ws.Range("G3").Formula = "=STDEV(C2:C11)"

# --- Column width / selection to match authoring state ---
# ColumnWidth is persisted with a constant +5px/MDW padding on save, so
# back the literal target width (21.5 chars) out of the desired value.
$ws.Columns.Item(5).ColumnWidth = 20.666666666666668
$ws.Range("E1:G3").Select()

# --- Page setup (portrait) ---
$ws.PageSetup.Orientation = 1

$wb.Save()
